$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G5: R15 resistor value is now a fixed literal instead of C5/10 ---
$ws.Range("G5").Value = 910

# --- Row 19: Div/(A*1000) -- switch from *100 to ROUND(*1000,0) ---
$ws.Range("A19").Value = "Div/(A*1000)"
$ws.Range("B19").Formula = "=ROUND(B15*1000,0)"
$ws.Range("C19").Formula = "=ROUND(C15*1000,0)"
$ws.Range("D19").Formula = "=ROUND(D15*1000,0)"
$ws.Range("F19").Formula = "=ROUND(F15*1000,0)"
$ws.Range("G19").Formula = "=ROUND(G15*1000,0)"
$ws.Range("H19").Formula = "=ROUND(H15*1000,0)"

# --- Row 20: Count -- new ADC codes ---
$ws.Range("C20").Value = 3541
$ws.Range("F20").Value = 60
$ws.Range("G20").Value = 3388
$ws.Range("H20").Value = 240

# --- Row 21: Count * 1000 (was Count * 100000) ---
$ws.Range("A21").Value = "Count * 1000"
$ws.Range("B21").Formula = "=B20*1000"
$ws.Range("C21").Formula = "=C20*1000"
$ws.Range("D21").Formula = "=D20*1000"

# --- Row 22: Count * 1000000 (new row label + formulas, replacing the old mA ratio) ---
$ws.Range("A22").Value = "Count * 1000000"
$ws.Range("B22").Formula = "=B20*1000000"
$ws.Range("C22").Formula = "=C20*1000000"
$ws.Range("D22").Formula = "=D20*1000000"
$ws.Range("F22").Formula = "=F20*1000000"
$ws.Range("G22").Formula = "=G20*1000000"
$ws.Range("H22").Formula = "=H20*1000000"

# --- Row 23: mA -- the final mA ratio (was row 22) ---
$ws.Range("A23").Value = "mA"
$ws.Range("B23").Formula = "=B22/B19"
$ws.Range("C23").Formula = "=C22/C19"
$ws.Range("D23").Formula = "=D22/D19"
$ws.Range("F23").Formula = "=F22/F19"
$ws.Range("G23").Formula = "=G22/G19"
$ws.Range("H23").Formula = "=H22/H19"

# --- Row 24/25: push the trailing blank formatting row down one ---
$ws.Range("J24").ClearContents()
$ws.Range("J25").Value = ""
